$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business
$ws.Range("H51").Value = 8876.6
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 8595.75
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 8595.75
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -9563.75

# Row 117: A Greater Grimoire
$ws.Range("H117").Value = 73871
$ws.Range("J117").Value = 73871
$ws.Range("L117").Value = 73871
$ws.Range("N117").Value = -83049

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 5689.795
$ws.Range("I132").Value = 4851.914
$ws.Range("K132").Value = 14555.742
$ws.Range("M132").Value = -12025.742

# Row 135: For Tired Minds
$ws.Range("H135").Value = 1174.5714
$ws.Range("I135").Value = 851.2222
$ws.Range("K135").Value = 7660.999800000001
$ws.Range("M135").Value = -5125.999800000001

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 5494.3335
$ws.Range("I141").Value = 3326.9048
$ws.Range("J141").Value = 20666.334
$ws.Range("K141").Value = 9980.714399999999
$ws.Range("L141").Value = 61999.00199999999
$ws.Range("M141").Value = -4800.714399999999
$ws.Range("N141").Value = -72359.00199999999

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate
$ws.Range("H4").Value = 328.5
$ws.Range("I4").Value = 423.33334
$ws.Range("J4").Value = 233.66667
$ws.Range("K4").Value = 423.33334
$ws.Range("L4").Value = 233.66667
$ws.Range("M4").Value = -307.33334
$ws.Range("N4").Value = -465.66667

# Row 44: Very Slow Array
$ws.Range("H44").Value = 71716
$ws.Range("J44").Value = 71716
$ws.Range("L44").Value = 71716
$ws.Range("N44").Value = -72692

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4954.5
$ws.Range("I61").Value = 3127.75
$ws.Range("J61").Value = 6781.25
$ws.Range("K61").Value = 3127.75
$ws.Range("L61").Value = 6781.25
$ws.Range("M61").Value = -2915.75
$ws.Range("N61").Value = -7205.25

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 3353.0715
$ws.Range("I110").Value = 1769.8422
$ws.Range("J110").Value = 6695.4443
$ws.Range("K110").Value = 1769.8422
$ws.Range("L110").Value = 6695.4443
$ws.Range("M110").Value = 275.1578
$ws.Range("N110").Value = -10785.4443

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3739.2273
$ws.Range("J132").Value = 7549.857
$ws.Range("L132").Value = 22649.571
$ws.Range("N132").Value = -27709.571

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4954.5
$ws.Range("I136").Value = 3127.75
$ws.Range("J136").Value = 6781.25
$ws.Range("K136").Value = 9383.25
$ws.Range("L136").Value = 20343.75
$ws.Range("M136").Value = -6833.25
$ws.Range("N136").Value = -25443.75

$ws = $wb.Worksheets.Item("BSM")
# Row 54: Get Me to the War on Time
$ws.Range("H54").Value = 54833.668
$ws.Range("I54").Value = 77500.5
$ws.Range("K54").Value = 77500.5
$ws.Range("M54").Value = -77016.5

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 21670336
$ws.Range("I105").Value = 5001602
$ws.Range("J105").Value = 25004084
$ws.Range("K105").Value = 5001602
$ws.Range("L105").Value = 25004084
$ws.Range("M105").Value = -4999855
$ws.Range("N105").Value = -25007578

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7790.636
$ws.Range("I134").Value = 7699.5
$ws.Range("J134").Value = 7900
$ws.Range("K134").Value = 23098.5
$ws.Range("L134").Value = 23700
$ws.Range("M134").Value = -20563.5
$ws.Range("N134").Value = -28770

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 4656.4546
$ws.Range("I7").Value = 6344.875
$ws.Range("K7").Value = 6344.875
$ws.Range("M7").Value = -6231.875

# Row 16: Raise the Roof
$ws.Range("H16").Value = 1524.1177
$ws.Range("I16").Value = 1613.7333
$ws.Range("J16").Value = 1453.3684
$ws.Range("K16").Value = 1613.7333
$ws.Range("L16").Value = 1453.3684
$ws.Range("M16").Value = -1326.7333
$ws.Range("N16").Value = -2027.3684

# Row 31: Wall Not Found
$ws.Range("H31").Value = 6540.724
$ws.Range("I31").Value = 5796.4546
$ws.Range("J31").Value = 6995.5557
$ws.Range("K31").Value = 5796.4546
$ws.Range("L31").Value = 6995.5557
$ws.Range("M31").Value = -5501.4546
$ws.Range("N31").Value = -7585.5557

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 6540.724
$ws.Range("I34").Value = 5796.4546
$ws.Range("J34").Value = 6995.5557
$ws.Range("K34").Value = 5796.4546
$ws.Range("L34").Value = 6995.5557
$ws.Range("M34").Value = -5594.4546
$ws.Range("N34").Value = -7399.5557

# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 31715.4
$ws.Range("J50").Value = 31715.4
$ws.Range("L50").Value = 31715.4
$ws.Range("N50").Value = -32965.4

# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 28215.2
$ws.Range("I51").Value = 19894
$ws.Range("J51").Value = 30295.5
$ws.Range("K51").Value = 19894
$ws.Range("L51").Value = 30295.5
$ws.Range("M51").Value = -19158
$ws.Range("N51").Value = -31767.5

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3798.875
$ws.Range("I58").Value = 4202.25
$ws.Range("J58").Value = 3664.4167
$ws.Range("K58").Value = 4202.25
$ws.Range("L58").Value = 3664.4167
$ws.Range("M58").Value = -3999.25
$ws.Range("N58").Value = -4070.4167

# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 28215.2
$ws.Range("I61").Value = 19894
$ws.Range("J61").Value = 30295.5
$ws.Range("K61").Value = 19894
$ws.Range("L61").Value = 30295.5
$ws.Range("M61").Value = -19546
$ws.Range("N61").Value = -30991.5

# Row 113: Patient Patients
$ws.Range("H113").Value = 1524.1177
$ws.Range("I113").Value = 1613.7333
$ws.Range("J113").Value = 1453.3684
$ws.Range("K113").Value = 1613.7333
$ws.Range("L113").Value = 1453.3684
$ws.Range("M113").Value = 556.2666999999999
$ws.Range("N113").Value = -5793.3684

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2858.543
$ws.Range("I132").Value = 2771.923
$ws.Range("J132").Value = 3108.7778
$ws.Range("K132").Value = 8315.769
$ws.Range("L132").Value = 9326.3334
$ws.Range("M132").Value = -5785.769
$ws.Range("N132").Value = -14386.3334

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2950.3845
$ws.Range("J134").Value = 3275.9
$ws.Range("L134").Value = 9827.700000000001
$ws.Range("N134").Value = -14897.7

# Row 136: Turali Quality
$ws.Range("H136").Value = 3798.875
$ws.Range("I136").Value = 4202.25
$ws.Range("J136").Value = 3664.4167
$ws.Range("K136").Value = 12606.75
$ws.Range("L136").Value = 10993.2501
$ws.Range("M136").Value = -10056.75
$ws.Range("N136").Value = -16093.2501

$ws = $wb.Worksheets.Item("CUL")
# Row 50: Moving Up in the World
$ws.Range("H50").Value = 1383.3334
$ws.Range("J50").Value = 1500
$ws.Range("L50").Value = 4500
$ws.Range("N50").Value = -5462

# Row 53: Rolanberry Fields Forever
$ws.Range("H53").Value = 1383.3334
$ws.Range("J53").Value = 1500
$ws.Range("L53").Value = 4500
$ws.Range("N53").Value = -5462

# Row 107: Slippery Service
$ws.Range("H107").Value = 340.8
$ws.Range("I107").Value = 245
$ws.Range("J107").Value = 364.75
$ws.Range("K107").Value = 735
$ws.Range("L107").Value = 1094.25
$ws.Range("M107").Value = 1185
$ws.Range("N107").Value = -4934.25

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 5997.1304
$ws.Range("I113").Value = 490.33334
$ws.Range("J113").Value = 7940.706
$ws.Range("K113").Value = 1471.00002
$ws.Range("L113").Value = 23822.118
$ws.Range("M113").Value = 698.9999800000001
$ws.Range("N113").Value = -28162.118

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 5478.483
$ws.Range("I122").Value = 4472.9565
$ws.Range("J122").Value = 9333
$ws.Range("K122").Value = 13418.8695
$ws.Range("L122").Value = 27999
$ws.Range("M122").Value = -10968.8695
$ws.Range("N122").Value = -32899

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 4838.875
$ws.Range("I126").Value = 2959.1428
$ws.Range("K126").Value = 8877.428400000001
$ws.Range("M126").Value = -6407.428400000001

# Row 132: On Board for Lar
$ws.Range("H132").Value = 6458.5625
$ws.Range("I132").Value = 1860.3334
$ws.Range("J132").Value = 12370.571
$ws.Range("K132").Value = 5581.0002
$ws.Range("L132").Value = 37111.713
$ws.Range("M132").Value = -3051.0002
$ws.Range("N132").Value = -42171.713

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 967
$ws.Range("I82").Value = 963.7
$ws.Range("K82").Value = 963.7
$ws.Range("M82").Value = -602.7

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 967
$ws.Range("I85").Value = 963.7
$ws.Range("K85").Value = 963.7
$ws.Range("M85").Value = 284.3

# Row 122: Hell on Leather
$ws.Range("H122").Value = 2293.8823
$ws.Range("I122").Value = 2192.6428
$ws.Range("J122").Value = 2766.3333
$ws.Range("K122").Value = 6577.928400000001
$ws.Range("L122").Value = 8298.999899999999
$ws.Range("M122").Value = -4127.928400000001
$ws.Range("N122").Value = -13198.9999

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 6773.522
$ws.Range("I132").Value = 2914.1428
$ws.Range("K132").Value = 8742.428400000001
$ws.Range("M132").Value = -6212.428400000001

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4310.3687
$ws.Range("I136").Value = 4136.4287
$ws.Range("J136").Value = 4797.4
$ws.Range("K136").Value = 12409.2861
$ws.Range("L136").Value = 14392.2
$ws.Range("M136").Value = -9859.286100000001
$ws.Range("N136").Value = -19492.2

$ws = $wb.Worksheets.Item("WVR")
# Row 56: Full Moon Fever
$ws.Range("H56").Value = 74104.5
$ws.Range("J56").Value = 74104.5
$ws.Range("L56").Value = 74104.5
$ws.Range("N56").Value = -75532.5

# Row 107: Flax Wax
$ws.Range("H107").Value = 868.5
$ws.Range("J107").Value = 1665
$ws.Range("L107").Value = 4995
$ws.Range("N107").Value = -8835

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2340.838
$ws.Range("I132").Value = 2263.9697
$ws.Range("J132").Value = 2975
$ws.Range("K132").Value = 6791.909100000001
$ws.Range("L132").Value = 8925
$ws.Range("M132").Value = -4261.909100000001
$ws.Range("N132").Value = -13985
